$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold + border + centered) from E1 to F1, then set header text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Populate time_taken values for each data row (F2:F232)
$timeTaken = @(
    "2021-10-05 13:42:32.301297",
    "2021-10-05 13:42:32.301308",
    "2021-10-05 13:42:32.301311",
    "2021-10-05 13:42:32.301314",
    "2021-10-05 13:42:32.301316",
    "2021-10-05 13:42:32.301319",
    "2021-10-05 13:42:32.301322",
    "2021-10-05 13:42:32.301324",
    "2021-10-05 13:42:32.301327",
    "2021-10-05 13:42:32.301329",
    "2021-10-05 13:42:32.301332",
    "2021-10-05 13:42:32.301334",
    "2021-10-05 13:42:32.301336",
    "2021-10-05 13:42:32.301339",
    "2021-10-05 13:42:32.301341",
    "2021-10-05 13:42:32.301344",
    "2021-10-05 13:42:32.301346",
    "2021-10-05 13:42:32.301349",
    "2021-10-05 13:42:32.301351",
    "2021-10-05 13:42:32.301354",
    "2021-10-05 13:42:32.301356",
    "2021-10-05 13:42:32.301359",
    "2021-10-05 13:42:32.301361",
    "2021-10-05 13:42:32.301364",
    "2021-10-05 13:42:32.301366",
    "2021-10-05 13:42:32.301369",
    "2021-10-05 13:42:32.301372",
    "2021-10-05 13:42:32.301374",
    "2021-10-05 13:42:32.301377",
    "2021-10-05 13:42:32.301379",
    "2021-10-05 13:42:32.301382",
    "2021-10-05 13:42:32.301384",
    "2021-10-05 13:42:32.301387",
    "2021-10-05 13:42:32.301390",
    "2021-10-05 13:42:32.301392",
    "2021-10-05 13:42:32.301395",
    "2021-10-05 13:42:32.301397",
    "2021-10-05 13:42:32.301400",
    "2021-10-05 13:42:32.301402",
    "2021-10-05 13:42:32.301405",
    "2021-10-05 13:42:32.301407",
    "2021-10-05 13:42:32.301410",
    "2021-10-05 13:42:32.301412",
    "2021-10-05 13:42:32.301415",
    "2021-10-05 13:42:32.301417",
    "2021-10-05 13:42:32.301420",
    "2021-10-05 13:42:32.301422",
    "2021-10-05 13:42:32.301424",
    "2021-10-05 13:42:32.301427",
    "2021-10-05 13:42:32.301429",
    "2021-10-05 13:42:32.301432",
    "2021-10-05 13:42:32.301434",
    "2021-10-05 13:42:32.301437",
    "2021-10-05 13:42:32.301440",
    "2021-10-05 13:42:32.301442",
    "2021-10-05 13:42:32.301444",
    "2021-10-05 13:42:32.301447",
    "2021-10-05 13:42:32.301449",
    "2021-10-05 13:42:32.301452",
    "2021-10-05 13:42:32.301454",
    "2021-10-05 13:42:32.301457",
    "2021-10-05 13:42:32.301459",
    "2021-10-05 13:42:32.301462",
    "2021-10-05 13:42:32.301464",
    "2021-10-05 13:42:32.301468",
    "2021-10-05 13:42:32.301471",
    "2021-10-05 13:42:32.301473",
    "2021-10-05 13:42:32.301475",
    "2021-10-05 13:42:32.301478",
    "2021-10-05 13:42:32.301480",
    "2021-10-05 13:42:32.301483",
    "2021-10-05 13:42:32.301485",
    "2021-10-05 13:42:32.301488",
    "2021-10-05 13:42:32.301490",
    "2021-10-05 13:42:32.301492",
    "2021-10-05 13:42:32.301495",
    "2021-10-05 13:42:32.301499",
    "2021-10-05 13:42:32.301502",
    "2021-10-05 13:42:32.301505",
    "2021-10-05 13:42:32.301507",
    "2021-10-05 13:42:32.301510",
    "2021-10-05 13:42:32.301512",
    "2021-10-05 13:42:32.301514",
    "2021-10-05 13:42:32.301517",
    "2021-10-05 13:42:32.301519",
    "2021-10-05 13:42:32.301522",
    "2021-10-05 13:42:32.301524",
    "2021-10-05 13:42:32.301527",
    "2021-10-05 13:42:32.301529",
    "2021-10-05 13:42:32.301532",
    "2021-10-05 13:42:32.301534",
    "2021-10-05 13:42:32.301536",
    "2021-10-05 13:42:32.301540",
    "2021-10-05 13:42:32.301543",
    "2021-10-05 13:42:32.301545",
    "2021-10-05 13:42:32.301548",
    "2021-10-05 13:42:32.301550",
    "2021-10-05 13:42:32.301553",
    "2021-10-05 13:42:32.301555",
    "2021-10-05 13:42:32.301558",
    "2021-10-05 13:42:32.301560",
    "2021-10-05 13:42:32.301562",
    "2021-10-05 13:42:32.301565",
    "2021-10-05 13:42:32.301567",
    "2021-10-05 13:42:32.301570",
    "2021-10-05 13:42:32.301572",
    "2021-10-05 13:42:32.301575",
    "2021-10-05 13:42:32.301577",
    "2021-10-05 13:42:32.301581",
    "2021-10-05 13:42:32.301584",
    "2021-10-05 13:42:32.301587",
    "2021-10-05 13:42:32.301589",
    "2021-10-05 13:42:32.301592",
    "2021-10-05 13:42:32.301594",
    "2021-10-05 13:42:32.301597",
    "2021-10-05 13:42:32.301599",
    "2021-10-05 13:42:32.301602",
    "2021-10-05 13:42:32.301604",
    "2021-10-05 13:42:32.301607",
    "2021-10-05 13:42:32.301609",
    "2021-10-05 13:42:32.301611",
    "2021-10-05 13:42:32.301614",
    "2021-10-05 13:42:32.301616",
    "2021-10-05 13:42:32.301619",
    "2021-10-05 13:42:32.301621",
    "2021-10-05 13:42:32.301624",
    "2021-10-05 13:42:32.301626",
    "2021-10-05 13:42:32.301628",
    "2021-10-05 13:42:32.301632",
    "2021-10-05 13:42:32.301635",
    "2021-10-05 13:42:32.301638",
    "2021-10-05 13:42:32.301640",
    "2021-10-05 13:42:32.301643",
    "2021-10-05 13:42:32.301645",
    "2021-10-05 13:42:32.301648",
    "2021-10-05 13:42:32.301650",
    "2021-10-05 13:42:32.301652",
    "2021-10-05 13:42:32.301655",
    "2021-10-05 13:42:32.301657",
    "2021-10-05 13:42:32.301660",
    "2021-10-05 13:42:32.301662",
    "2021-10-05 13:42:32.301664",
    "2021-10-05 13:42:32.301667",
    "2021-10-05 13:42:32.301669",
    "2021-10-05 13:42:32.301672",
    "2021-10-05 13:42:32.301674",
    "2021-10-05 13:42:32.301677",
    "2021-10-05 13:42:32.301679",
    "2021-10-05 13:42:32.301682",
    "2021-10-05 13:42:32.301685",
    "2021-10-05 13:42:32.301687",
    "2021-10-05 13:42:32.301690",
    "2021-10-05 13:42:32.301692",
    "2021-10-05 13:42:32.301695",
    "2021-10-05 13:42:32.301697",
    "2021-10-05 13:42:32.301700",
    "2021-10-05 13:42:32.301702",
    "2021-10-05 13:42:32.301704",
    "2021-10-05 13:42:32.301707",
    "2021-10-05 13:42:32.301709",
    "2021-10-05 13:42:32.301712",
    "2021-10-05 13:42:32.301714",
    "2021-10-05 13:42:32.301717",
    "2021-10-05 13:42:32.301719",
    "2021-10-05 13:42:32.301721",
    "2021-10-05 13:42:32.301724",
    "2021-10-05 13:42:32.301726",
    "2021-10-05 13:42:32.301729",
    "2021-10-05 13:42:32.301731",
    "2021-10-05 13:42:32.301734",
    "2021-10-05 13:42:32.301736",
    "2021-10-05 13:42:32.301738",
    "2021-10-05 13:42:32.301742",
    "2021-10-05 13:42:32.301745",
    "2021-10-05 13:42:32.301748",
    "2021-10-05 13:42:32.301750",
    "2021-10-05 13:42:32.301753",
    "2021-10-05 13:42:32.301755",
    "2021-10-05 13:42:32.301757",
    "2021-10-05 13:42:32.301760",
    "2021-10-05 13:42:32.301762",
    "2021-10-05 13:42:32.301765",
    "2021-10-05 13:42:32.301767",
    "2021-10-05 13:42:32.301816",
    "2021-10-05 13:42:32.301828",
    "2021-10-05 13:42:32.301832",
    "2021-10-05 13:42:32.301834",
    "2021-10-05 13:42:32.301837",
    "2021-10-05 13:42:32.301840",
    "2021-10-05 13:42:32.301843",
    "2021-10-05 13:42:32.301845",
    "2021-10-05 13:42:32.301848",
    "2021-10-05 13:42:32.301851",
    "2021-10-05 13:42:32.301853",
    "2021-10-05 13:42:32.301856",
    "2021-10-05 13:42:32.301858",
    "2021-10-05 13:42:32.301861",
    "2021-10-05 13:42:32.301863",
    "2021-10-05 13:42:32.301866",
    "2021-10-05 13:42:32.301869",
    "2021-10-05 13:42:32.301872",
    "2021-10-05 13:42:32.301875",
    "2021-10-05 13:42:32.301878",
    "2021-10-05 13:42:32.301881",
    "2021-10-05 13:42:32.301883",
    "2021-10-05 13:42:32.301886",
    "2021-10-05 13:42:32.301888",
    "2021-10-05 13:42:32.301891",
    "2021-10-05 13:42:32.301893",
    "2021-10-05 13:42:32.301896",
    "2021-10-05 13:42:32.301898",
    "2021-10-05 13:42:32.301901",
    "2021-10-05 13:42:32.301903",
    "2021-10-05 13:42:32.301906",
    "2021-10-05 13:42:32.301908",
    "2021-10-05 13:42:32.301911",
    "2021-10-05 13:42:32.301914",
    "2021-10-05 13:42:32.301916",
    "2021-10-05 13:42:32.301919",
    "2021-10-05 13:42:32.301921",
    "2021-10-05 13:42:32.301924",
    "2021-10-05 13:42:32.301926",
    "2021-10-05 13:42:32.301929",
    "2021-10-05 13:42:32.301931",
    "2021-10-05 13:42:32.301933",
    "2021-10-05 13:42:32.301936",
    "2021-10-05 13:42:32.301938",
    "2021-10-05 13:42:32.301941",
    "2021-10-05 13:42:32.301943",
    "2021-10-05 13:42:32.301946",
    "2021-10-05 13:42:32.301948"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
